# Servery_obhajoba.pptx — "Add files via upload"
#
# Slide 3  ("VYBER TEMY" title box):     enlarge the accent letter "Y-acute"
#                                        and shrink the (spAutoFit) box height
#                                        to match; nudge the rotated corner-
#                                        flourish group that sits above it.
# Slide 9  ("DAKUJEM ZA POZORNOST" box): enlarge the leading "D-caron" and
#                                        trailing "T-caron" accent letters and
#                                        shrink the (spAutoFit) box height to
#                                        match.
#
# EMU <-> point conversion is the usual 12700 EMU/pt; the handful of literal
# point constants below carry a few extra decimal digits so that the COM
# host's internal float32 Left/Top/Height storage round-trips back to the
# *exact* target EMU value instead of landing 1 EMU off by truncation.

$p  = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 - "VYBER TEMY"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# TextBox 10 holds the title text.
$title = $s3.Shapes.Item(7)
$titleRange = $title.TextFrame.TextRange

# Split out the 2nd character ("Y-acute") into its own bigger, bold run.
$accent = $titleRange.Characters(2, 1)
$accent.Font.Size = 110
$accent.Font.Bold = $true

# The textbox auto-fits its height to the text; the bigger accent letter no
# longer drives the overall box height once it matches the sibling slide's
# identical title box, so pin the height to the author's resized value.
$title.Height = 105.05598425196851

# Group 12 is the rotated corner-flourish decoration; it shifts slightly to
# re-align now that the title box above it is shorter.
$corner = $s3.Shapes.Item(9)
$corner.Left = -115.95512011023621
$corner.Top = -228.1708661417323

# ---------------------------------------------------------------------
# Slide 9 - "DAKUJEM ZA POZORNOST"
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)

$thanks = $s9.Shapes.Item(1)
$thanksRange = $thanks.TextFrame.TextRange

# Split out the first character ("D-caron") ...
$first = $thanksRange.Characters(1, 1)
$first.Font.Size = 134
$first.Font.Bold = $true

# ... and the last character ("T-caron") into their own bigger, bold runs.
$last = $thanksRange.Characters(20, 1)
$last.Font.Size = 134
$last.Font.Bold = $true

# Pin the autosized box height to the author's resized value.
$thanks.Height = 250.42188976377952
